# Refresh the crypto price/volume table (columns D = Price, E = Volume(1h))
# with the latest scraped figures. Values that look like plain decimal
# numbers are prefixed with a leading apostrophe so Excel stores them as
# text (preserving formatting such as trailing zeros) instead of coercing
# them to numbers, matching the original inline-string cell types.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.993.64"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "1.641.36"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'213.06"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").Value = "'0.524"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'23.62"
$ws.Range("E8").Value = "  +1.76%  "
$ws.Range("D9").Value = "'0.259"
$ws.Range("E9").Value = "  -1.85%  "
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").Value = "'0.0882"
$ws.Range("E11").Value = "  +2.41%  "
$ws.Range("D12").Value = "1.873.88"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").Value = "1.642.01"
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("E14").Value = "  +3.79%  "
$ws.Range("D15").Value = "'4.09"
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("D16").Value = "'65.90"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("D17").Value = "27.990.30"
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("D18").Value = "'235.74"
$ws.Range("E18").Value = "  +2.50%  "
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "'10.68"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").Value = "'4.38"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("E24").Value = "  -2.21%  "
$ws.Range("D25").Value = "'151.16"
$ws.Range("D26").Value = "'6.96"
$ws.Range("E26").Value = "  +1.40%  "
$ws.Range("D27").Value = "'15.70"
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("D34").Value = "1.418.85"
$ws.Range("E34").Value = "  -3.68%  "
$ws.Range("E35").Value = "  +2.34%  "
$ws.Range("D36").Value = "'2.35"
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("E37").Value = "  +1.58%  "
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("D39").Value = "'0.559"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "'0.902"
$ws.Range("E40").Value = "  -4.88%  "
$ws.Range("E41").Value = "  +1.13%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  +6.69%  "
$ws.Range("D44").Value = "'66.65"
$ws.Range("E44").Value = "  -1.91%  "
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").Value = "1.782.96"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("D48").Value = "'87.87"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").Value = "'7.62"
$ws.Range("E51").Value = "  -1.38%  "
